$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Scroll view back to top-left (removes topLeftCell="AB1") ---
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1

# --- Row 2 ---
$ws.Range("D2").NumberFormat = "0.00"
$ws.Range("R2").Value = 2
$ws.Range("AL2").Value = 35296

# --- Row 3 ---
$ws.Range("D3").NumberFormat = "0.00"
$ws.Range("Q3").Value = 3
$ws.Range("R3").ClearContents()
$ws.Range("AM3").ClearContents()

# --- Row 4 ---
$ws.Range("D4").NumberFormat = "0.00"
$ws.Range("Q4").Value = 1
$ws.Range("R4").ClearContents()

# --- Row 5 ---
$ws.Range("D5").NumberFormat = "0.000"
$ws.Range("D5").Value = 44.999000000000002
$ws.Range("Q5").Value = 2
$ws.Range("AL5").Value = 38221
$ws.Range("AM5").ClearContents()

# --- Row 6 ---
$ws.Range("D6").NumberFormat = "0.00"
$ws.Range("D6").Value = 45
$ws.Range("R6").ClearContents()

# --- Row 7 ---
$ws.Range("D7").NumberFormat = "0.00"
$ws.Range("D7").Value = 55.01
$ws.Range("Q7").Value = 2
$ws.Range("R7").ClearContents()
$ws.Range("AL7").Value = 36769

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
